$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted into the daily log.
# Insert a new row at 61, shifting all existing rows 61..126 down to 62..127.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new observation.
$ws.Cells.Item(61, 1).Value = 9
$ws.Cells.Item(61, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(61, 3).Value = "Metropolitana"
$ws.Cells.Item(61, 4).Value = 44790
$ws.Cells.Item(61, 5).Value = 13
$ws.Cells.Item(61, 6).Value = 100112022
$ws.Cells.Item(61, 7).Value = "Arveja Verde"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 36
$ws.Cells.Item(61, 11).Value = 41000
$ws.Cells.Item(61, 12).Value = 41000
$ws.Cells.Item(61, 13).Value = 41000
$ws.Cells.Item(61, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(61, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(61, 16).Value = 1640
$ws.Cells.Item(61, 17).Value = 25
$ws.Cells.Item(61, 18).Value = "Hortaliza"
